# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to reflect the latest generated numbers (commit: Update gh-pages to
# output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row => new value for column F
$updates = @{
    3  = 1361
    4  = 155
    6  = 232
    8  = 15
    11 = 4570
    12 = 6843
    16 = 571
    18 = 4133
    19 = 522
    20 = 75
    21 = 59
    22 = 2714
    24 = 549
    25 = 169
    26 = 356
    27 = 367
    29 = 225
    30 = 39
    31 = 1626
    34 = 161
    35 = 82
    39 = 90
    40 = 113
    41 = 644
    42 = 13
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
